# Automatically write the input values of vin and id into the excel worksheets.

$wb = $excel.ActiveWorkbook

# --- Cars sheet: append new rows of vin values in column A ---
$carsWs = $wb.Worksheets.Item("Cars")
$vins = @(
    "25F9N5NXM4IGC4342",
    "LHMONS7M2F1OKPRT3",
    "DTWGB2ODEO4XOALNU",
    "76178I4Z4JVTBKD0R",
    "BPEVHNWEWII29IHOX",
    "AAV7P62EF2GFVF1XL",
    "SQF6J5UHM1E"
)
$row = 2
foreach ($vin in $vins) {
    $carsWs.Range("A" + $row).Value = $vin
    $row = $row + 1
}

# --- Sellers sheet: append new seller id value in column A ---
$sellersWs = $wb.Worksheets.Item("Sellers")
$sellersWs.Range("A2").Value = "DJ41R7T8LTT"
$sellersWs.Columns.Item(1).ColumnWidth = 29.5

# --- Sellers becomes the active/selected sheet ---
$sellersWs.Activate()
